$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 1.98
$ws.Range("H4").Value = 3.35
$ws.Range("I4").Value = 3.55
$ws.Range("L4").Value = 1.36
$ws.Range("M4").Value = 2.65
$ws.Range("N4").Value = 2.05
$ws.Range("O4").Value = 1.6
$ws.Range("P4").Value = 1.45
$ws.Range("Q4").Value = 2.37
$ws.Range("R4").Value = 1.9
$ws.Range("S4").Value = 1.72
$ws.Range("T4").Value = 6.3
$ws.Range("U4").Value = 8.5
$ws.Range("V4").Value = 8.75
$ws.Range("W4").Value = 17
$ws.Range("X4").Value = 17.5
$ws.Range("Y4").Value = 32
$ws.Range("Z4").Value = 8.25
$ws.Range("AA4").Value = 6.5
$ws.Range("AC4").Value = 90
$ws.Range("AD4").Value = 9
$ws.Range("AE4").Value = 17.5
$ws.Range("AF4").Value = 12.5
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 37
$ws.Range("AI4").Value = 50
$ws.Range("AJ4").Value = 900

# Row 6
$ws.Range("G6").Value = 1.42
$ws.Range("H6").Value = 4.5
$ws.Range("I6").Value = 5.75
$ws.Range("J6").Value = 1.06
$ws.Range("K6").Value = 10
$ws.Range("P6").Value = 1.37
$ws.Range("R6").Value = 2.25
$ws.Range("S6").Value = 1.57
$ws.Range("T6").Value = 5.5
$ws.Range("U6").Value = 6
$ws.Range("V6").Value = 9
$ws.Range("W6").Value = 9
$ws.Range("AA6").Value = 9.5
$ws.Range("AB6").Value = 26
$ws.Range("AD6").Value = 13
$ws.Range("AE6").Value = 34
$ws.Range("AF6").Value = 19
$ws.Range("AG6").Value = 81

# Row 7
$ws.Range("G7").Value = 3.8
$ws.Range("I7").Value = 1.85
$ws.Range("J7").Value = 1.05
$ws.Range("K7").Value = 11
$ws.Range("AA7").Value = 7
$ws.Range("AG7").Value = 15
$ws.Range("AI7").Value = 26
$ws.Range("AJ7").Value = 251

# Row 8
$ws.Range("L8").Value = 1.45
$ws.Range("M8").Value = 2.37
$ws.Range("N8").Value = 2.32
$ws.Range("P8").Value = 1.47
$ws.Range("Q8").Value = 2.32
$ws.Range("R8").Value = 2.07
$ws.Range("S8").Value = 1.6
$ws.Range("T8").Value = 5.4
$ws.Range("U8").Value = 7.7
$ws.Range("V8").Value = 9
$ws.Range("W8").Value = 16
$ws.Range("X8").Value = 18.5
$ws.Range("Y8").Value = 40
$ws.Range("Z8").Value = 6.9
$ws.Range("AC8").Value = 120

# Row 10
$ws.Range("H10").Value = 5.4
$ws.Range("I10").Value = 13.5
$ws.Range("K10").Value = 8.25
$ws.Range("L10").Value = 1.23
$ws.Range("M10").Value = 3.8
$ws.Range("N10").Value = 1.7
$ws.Range("O10").Value = 2.05
$ws.Range("P10").Value = 1.35
$ws.Range("Q10").Value = 2.95
$ws.Range("R10").Value = 2.4
$ws.Range("S10").Value = 1.5
$ws.Range("T10").Value = 6.3
$ws.Range("X10").Value = 11.5
$ws.Range("Y10").Value = 40
$ws.Range("Z10").Value = 8.25
$ws.Range("AA10").Value = 11.25
$ws.Range("AC10").Value = 200
$ws.Range("AD10").Value = 28
$ws.Range("AE10").Value = 110
$ws.Range("AF10").Value = 45
$ws.Range("AG10").Value = 600
$ws.Range("AH10").Value = 250
$ws.Range("AI10").Value = 200

# Row 12
$ws.Range("G12").Value = 4.65
$ws.Range("H12").Value = 3.5
$ws.Range("I12").Value = 1.72
$ws.Range("M12").Value = 3.25
$ws.Range("N12").Value = 1.87
$ws.Range("O12").Value = 1.83
$ws.Range("P12").Value = 1.42
$ws.Range("Q12").Value = 2.65
$ws.Range("T12").Value = 12.5
$ws.Range("U12").Value = 27
$ws.Range("V12").Value = 15
$ws.Range("W12").Value = 80
$ws.Range("X12").Value = 45
$ws.Range("Y12").Value = 50
$ws.Range("AA12").Value = 6.7
$ws.Range("AB12").Value = 15
$ws.Range("AD12").Value = 6.9
$ws.Range("AF12").Value = 8
$ws.Range("AG12").Value = 14
$ws.Range("AH12").Value = 13.5
$ws.Range("AI12").Value = 25
$ws.Range("AJ12").Value = 600

# Row 15
$ws.Range("G15").Value = 1.18
$ws.Range("H15").Value = 5.2
$ws.Range("V15").Value = 8
$ws.Range("W15").Value = 5.6
$ws.Range("Z15").Value = 12
$ws.Range("AB15").Value = 24
$ws.Range("AD15").Value = 28
$ws.Range("AE15").Value = 100
$ws.Range("AI15").Value = 120

# Row 17
$ws.Range("T17").Value = 6.5
$ws.Range("U17").Value = 9
$ws.Range("W17").Value = 18
$ws.Range("X17").Value = 17.5
$ws.Range("Y17").Value = 30
$ws.Range("AB17").Value = 14.5
$ws.Range("AC17").Value = 70
$ws.Range("AD17").Value = 10.25
$ws.Range("AE17").Value = 21
$ws.Range("AH17").Value = 35
$ws.Range("AI17").Value = 40
$ws.Range("AJ17").Value = 600

# Row 19
$ws.Range("L19").Value = 1.29
$ws.Range("M19").Value = 3.5
$ws.Range("N19").Value = 1.9
$ws.Range("O19").Value = 1.9

# Row 21
$ws.Range("G21").Value = 3.55
$ws.Range("H21").Value = 4
$ws.Range("I21").Value = 1.8
$ws.Range("N21").Value = 1.34
$ws.Range("O21").Value = 2.72
$ws.Range("R21").Value = 1.36
$ws.Range("S21").Value = 2.67
$ws.Range("T21").Value = 20
$ws.Range("U21").Value = 27
$ws.Range("V21").Value = 13
$ws.Range("W21").Value = 55
$ws.Range("X21").Value = 26
$ws.Range("Y21").Value = 23
$ws.Range("Z21").Value = 22
$ws.Range("AA21").Value = 9.25
$ws.Range("AB21").Value = 11.25
$ws.Range("AC21").Value = 29
$ws.Range("AD21").Value = 14
$ws.Range("AE21").Value = 13.5
$ws.Range("AF21").Value = 9
$ws.Range("AG21").Value = 18.5
$ws.Range("AH21").Value = 12.5
$ws.Range("AI21").Value = 15.5
$ws.Range("AJ21").Value = 120

# Row 24
$ws.Range("J24").Value = 1.05
$ws.Range("L24").Value = 1.33
$ws.Range("R24").Value = 1.87
$ws.Range("S24").Value = 1.87

# Row 30
$ws.Range("J30").Value = 1.08
$ws.Range("K30").Value = 8
$ws.Range("N30").Value = 2.25
$ws.Range("O30").Value = 1.62
$ws.Range("P30").Value = 1.5

# Row 31
$ws.Range("R31").Value = 1.8
$ws.Range("S31").Value = 1.8

# Row 32
$ws.Range("P32").Value = 1.5

# Row 34
$ws.Range("P34").Value = 1.29
$ws.Range("R34").Value = 1.5
$ws.Range("S34").Value = 2.37

# Row 36
$ws.Range("P36").Value = 1.3

# Row 37
$ws.Range("R37").Value = 1.63

# Row 38
$ws.Range("J38").Value = 1.03
$ws.Range("L38").Value = 1.22
$ws.Range("N38").Value = 1.77
$ws.Range("O38").Value = 1.92

# Row 39
$ws.Range("J39").Value = 1.03
$ws.Range("L39").Value = 1.19
$ws.Range("R39").Value = 1.67

# Row 42
$ws.Range("J42").Value = 1.05
$ws.Range("K42").Value = 11

# Row 44
$ws.Range("C44").Value = "23:15"
$ws.Range("G44").Value = 1.91
$ws.Range("J44").Value = 1.02
$ws.Range("L44").Value = 1.13
$ws.Range("M44").Value = 6
$ws.Range("R44").Value = 1.41
$ws.Range("S44").Value = 2.62

# Row 46
$ws.Range("G46").Value = 2.3
$ws.Range("J46").Value = 1.05
$ws.Range("L46").Value = 1.25

# Row 48
$ws.Range("G48").Value = 2.18
$ws.Range("H48").Value = 3.3
$ws.Range("I48").Value = 2.95
$ws.Range("K48").Value = 7
$ws.Range("P48").Value = 1.44
$ws.Range("R48").Value = 1.82
$ws.Range("T48").Value = 7.3
$ws.Range("Z48").Value = 7
$ws.Range("AC48").Value = 75
$ws.Range("AE48").Value = 15
$ws.Range("AG48").Value = 37

# Row 49
$ws.Range("G49").Value = 2.12
$ws.Range("H49").Value = 3.4
$ws.Range("I49").Value = 2.95
$ws.Range("N49").Value = 1.91
$ws.Range("O49").Value = 1.8
$ws.Range("P49").Value = 1.4
$ws.Range("Q49").Value = 2.7
$ws.Range("T49").Value = 7.6
$ws.Range("U49").Value = 10.25
$ws.Range("V49").Value = 9
$ws.Range("W49").Value = 20
$ws.Range("X49").Value = 17.5
$ws.Range("Y49").Value = 29
$ws.Range("AA49").Value = 6.7
$ws.Range("AB49").Value = 15
$ws.Range("AD49").Value = 9.25
$ws.Range("AE49").Value = 15
$ws.Range("AF49").Value = 11
$ws.Range("AG49").Value = 37
$ws.Range("AH49").Value = 26
$ws.Range("AI49").Value = 35
$ws.Range("AJ49").Value = 600
